$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.955.51"
$ws.Range("E2").Value = "  -1.12%  "
$ws.Range("D3").Value = "2.663.43"
$ws.Range("E3").Value = "  +1.20%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.39%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.524"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.53%  "
$ws.Range("D9").Value = "2.663.34"
$ws.Range("E9").Value = "  +1.27%  "
$ws.Range("E10").Value = "  -1.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.170"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.356"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.00"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.12%  "
$ws.Range("D14").Value = "3.149.41"
$ws.Range("E14").Value = "  +0.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000185"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.24%  "
$ws.Range("D16").Value = "71.731.39"
$ws.Range("E16").Value = "  -1.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.50%  "
$ws.Range("D18").Value = "2.662.59"
$ws.Range("E18").Value = "  +0.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "371.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.16"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.35%  "
$ws.Range("E23").Value = "  +1.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "72.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.46%  "
$ws.Range("D28").Value = "2.799.68"
$ws.Range("E28").Value = "  +1.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.51%  "
$ws.Range("D30").Value = "0.0₃0970"
$ws.Range("E30").Value = "  +0.96%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.06"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "500.50"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.03%  "
$ws.Range("E33").Value = "  -2.53%  "
$ws.Range("E34").Value = "  -0.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "163.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.51"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.05"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("E39").Value = "  -2.19%  "
$ws.Range("E40").Value = "  -2.16%  "
$ws.Range("E41").Value = "  -3.29%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.99"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.62%  "
$ws.Range("E44").Value = "  -2.05%  "
$ws.Range("E45").Value = "  -0.10%  "
$ws.Range("E46").Value = "  -0.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "156.03"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.559"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.73"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.27%  "
$ws.Range("E50").Value = "  +1.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0753"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.59%  "
